{"js": "const oldText =\n  \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od Souhv\u011bzd\u00ed Perseus 2022: 16. \u2013 25. ledna, 7. \u2013 16. listopadu, 6. \u2013 15. prosince\";\nconst newText =\n  \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od 16. \u2013 25. ledna, 7. \u2013 16. listopadu, 6. \u2013 15. prosince. P\u0159i pozorov\u00e1n\u00ed pou\u017eijte hv\u011bzdy oblohy, kter\u00e9 zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed Perseus.16. \u2013 25. ledna, 7. \u2013 16. listopadu, 6. \u2013 15. prosince\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$old = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od Souhv\u011bzd\u00ed Perseus 2022: 16. \u2013 25. ledna, 7. \u2013 16. listopadu, 6. \u2013 15. prosince\"\n$new = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od 16. \u2013 25. ledna, 7. \u2013 16. listopadu, 6. \u2013 15. prosince. P\u0159i pozorov\u00e1n\u00ed pou\u017eijte hv\u011bzdy oblohy, kter\u00e9 zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed Perseus.16. \u2013 25. ledna, 7. \u2013 16. listopadu, 6. \u2013 15. prosince\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $old\n$find.Forward = $true\n$find.Wrap = 0\n\nwhile ($find.Execute()) {\n    $rng.Text = $new\n    $rng.Collapse(0)\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Forward = $true\n    $find.Wrap = 0\n}\n"}
